$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.561.31"
$ws.Range("E2").Value = "  -2.84%  "

$ws.Range("D3").Value = "3.724.47"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.01"
$ws.Range("E5").Value = "  -3.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.38"
$ws.Range("E6").Value = "  -3.79%  "

$ws.Range("D7").Value = "3.738.10"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  -2.02%  "

$ws.Range("E10").Value = "  -5.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("E11").Value = "  -6.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  -4.74%  "

$ws.Range("E13").Value = "  -6.47%  "

$ws.Range("E14").Value = "  -5.01%  "

$ws.Range("D15").Value = "4.343.18"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").Value = "3.716.85"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("D17").Value = "67.499.50"
$ws.Range("E17").Value = "  -3.02%  "

$ws.Range("E18").Value = "  -5.16%  "

$ws.Range("E19").Value = "  -4.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.07"
$ws.Range("E20").Value = "  -1.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.25"
$ws.Range("E21").Value = "  -2.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.95"
$ws.Range("E22").Value = "  -2.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.711"
$ws.Range("E23").Value = "  -1.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.32"
$ws.Range("E24").Value = "  -2.92%  "

$ws.Range("E25").Value = "  -9.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000140"
$ws.Range("E26").Value = "  +4.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.09"
$ws.Range("E27").Value = "  -5.93%  "

$ws.Range("E28").Value = "  -7.97%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("E30").Value = "  +0.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.14"
$ws.Range("E31").Value = "  +5.86%  "

$ws.Range("E32").Value = "  -4.00%  "

$ws.Range("E33").Value = "  -4.38%  "

$ws.Range("E34").Value = "  -4.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.992"
$ws.Range("E36").Value = "  -5.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.134"
$ws.Range("E37").Value = "  -2.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.67"
$ws.Range("E38").Value = "  -6.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.322"
$ws.Range("E39").Value = "  -7.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "447.15"
$ws.Range("E40").Value = "  +0.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "48.82"
$ws.Range("E41").Value = "  -1.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.98"
$ws.Range("E42").Value = "  -3.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.85"
$ws.Range("E43").Value = "  -6.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.22"
$ws.Range("E44").Value = "  -3.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.07"
$ws.Range("E45").Value = "  -7.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "140.69"
$ws.Range("E46").Value = "  +1.18%  "

$ws.Range("D47").Value = "2.782.62"
$ws.Range("E47").Value = "  -5.68%  "

$ws.Range("E49").Value = "  -3.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.55"
$ws.Range("E50").Value = "  -5.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.10"
$ws.Range("E51").Value = "  +7.80%  "
